$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 779
$ws.Range("F3").Value = 2816
$ws.Range("F5").Value = 1962
$ws.Range("F7").Value = 50
$ws.Range("F9").Value = 285
$ws.Range("F11").Value = 11746
$ws.Range("F12").Value = 6674
$ws.Range("F19").Value = 927
$ws.Range("F21").Value = 279
$ws.Range("F22").Value = 932
$ws.Range("F26").Value = 499
$ws.Range("F30").Value = 231
$ws.Range("F31").Value = 268
$ws.Range("F32").Value = 310
$ws.Range("F33").Value = 5034
$ws.Range("F35").Value = 1247
$ws.Range("F36").Value = 240
$ws.Range("F37").Value = 560
$ws.Range("F38").Value = 207
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 3690
$ws.Range("F23").Value = 15
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9066
$ws.Range("F3").Value = 510
$ws.Range("F4").Value = 1837
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 510
$ws.Range("F3").Value = 1837
$ws.Range("F4").Value = 779
$ws.Range("F5").Value = 2816
$ws.Range("F12").Value = 50
$ws.Range("F15").Value = 285
$ws.Range("F17").Value = 11746
$ws.Range("F18").Value = 3690
$ws.Range("F19").Value = 6674
$ws.Range("F28").Value = 279
$ws.Range("F29").Value = 932
$ws.Range("F35").Value = 231
$ws.Range("F36").Value = 268
$ws.Range("F40").Value = 1247
$ws.Range("F41").Value = 240
$ws.Range("F43").Value = 207
$ws.Range("F48").Value = 15
